$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Dhruv Goyani's row (row 9): Knowledge Area & Skill selection changed
$ws.Range("B9").Value = "Computing Foundations"
$ws.Range("C9").Value = "Algorithms and Complexity"

# Give every data row (1-32) an explicit 15.75pt custom height
for ($r = 1; $r -le 32; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# Move the active selection to C10
$ws.Range("C10").Select()
